$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old sample data (A1:B20) from the sheet.
$ws.Cells.Clear()

# Leave only the finished note in F6.
$ws.Range("F6").Value = "ide to"
